$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 804.5
$ws.Range("J129").Value = 944.625
$ws.Range("L129").Value = 2833.875
$ws.Range("N129").Value = -12833.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3384.9167
$ws.Range("I138").Value = 2851.2058
$ws.Range("J138").Value = 4082.8462
$ws.Range("K138").Value = 8553.617400000001
$ws.Range("L138").Value = 12248.5386
$ws.Range("M138").Value = -3413.617400000001
$ws.Range("N138").Value = -22528.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1119.7273
$ws.Range("I45").Value = 892.6667
$ws.Range("J45").Value = 1276.9231
$ws.Range("K45").Value = 892.6667
$ws.Range("L45").Value = 1276.9231
$ws.Range("M45").Value = -515.6667
$ws.Range("N45").Value = -2030.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3475637.8
$ws.Range("I61").Value = 5293858
$ws.Range("J61").Value = 4490.364
$ws.Range("K61").Value = 5293858
$ws.Range("L61").Value = 4490.364
$ws.Range("M61").Value = -5293646
$ws.Range("N61").Value = -4914.364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 670.29266
$ws.Range("I97").Value = 564.32355
$ws.Range("J97").Value = 1185
$ws.Range("K97").Value = 564.32355
$ws.Range("L97").Value = 1185
$ws.Range("M97").Value = -68.32354999999995
$ws.Range("N97").Value = -2177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2192.4285
$ws.Range("I102").Value = 2252.5386
$ws.Range("K102").Value = 2252.5386
$ws.Range("M102").Value = -630.5385999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1227.1765
$ws.Range("I110").Value = 1100.1538
$ws.Range("J110").Value = 1640
$ws.Range("K110").Value = 1100.1538
$ws.Range("L110").Value = 1640
$ws.Range("M110").Value = 944.8462
$ws.Range("N110").Value = -5730

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1381933.8
$ws.Range("I132").Value = 2030567
$ws.Range("J132").Value = 3588.25
$ws.Range("K132").Value = 6091701
$ws.Range("L132").Value = 10764.75
$ws.Range("M132").Value = -6089171
$ws.Range("N132").Value = -15824.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3475637.8
$ws.Range("I136").Value = 5293858
$ws.Range("J136").Value = 4490.364
$ws.Range("K136").Value = 15881574
$ws.Range("L136").Value = 13471.092
$ws.Range("M136").Value = -15879024
$ws.Range("N136").Value = -18571.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 302.41666
$ws.Range("I94").Value = 285.55554
$ws.Range("J94").Value = 353
$ws.Range("K94").Value = 285.55554
$ws.Range("L94").Value = 353
$ws.Range("M94").Value = 165.44446
$ws.Range("N94").Value = -1255

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1480.4
$ws.Range("I99").Value = 1409.8182
$ws.Range("J99").Value = 1674.5
$ws.Range("K99").Value = 1409.8182
$ws.Range("L99").Value = 1674.5
$ws.Range("M99").Value = 88.18180000000007
$ws.Range("N99").Value = -4670.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1726.1538
$ws.Range("I105").Value = 1493.3334
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 1493.3334
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = 253.6666
$ws.Range("N105").Value = -5744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 347241.72
$ws.Range("I107").Value = 544873.4399999999
$ws.Range("J107").Value = 4680
$ws.Range("K107").Value = 544873.4399999999
$ws.Range("L107").Value = 4680
$ws.Range("M107").Value = -542953.4399999999
$ws.Range("N107").Value = -8520

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15398340
$ws.Range("I134").Value = 16681285
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 50043855
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -50041320
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9488.556
$ws.Range("J51").Value = 9685.286
$ws.Range("L51").Value = 9685.286
$ws.Range("N51").Value = -11157.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2988.2742
$ws.Range("I58").Value = 1239.9375
$ws.Range("J58").Value = 3596.3914
$ws.Range("K58").Value = 1239.9375
$ws.Range("L58").Value = 3596.3914
$ws.Range("M58").Value = -1036.9375
$ws.Range("N58").Value = -4002.3914

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 10201.714
$ws.Range("J60").Value = 10201.714
$ws.Range("L60").Value = 10201.714
$ws.Range("N60").Value = -11223.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9488.556
$ws.Range("J61").Value = 9685.286
$ws.Range("L61").Value = 9685.286
$ws.Range("N61").Value = -10381.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2988.2742
$ws.Range("I136").Value = 1239.9375
$ws.Range("J136").Value = 3596.3914
$ws.Range("K136").Value = 3719.8125
$ws.Range("L136").Value = 10789.1742
$ws.Range("M136").Value = -1169.8125
$ws.Range("N136").Value = -15889.1742

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7247535
$ws.Range("J122").Value = 1779.7727
$ws.Range("L122").Value = 16017.9543
$ws.Range("N122").Value = -20917.9543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2954.745
$ws.Range("I131").Value = 20274.334
$ws.Range("J131").Value = 1872.2709
$ws.Range("K131").Value = 60823.00199999999
$ws.Range("L131").Value = 5616.8127
$ws.Range("M131").Value = -55783.00199999999
$ws.Range("N131").Value = -15696.8127

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4704.303
$ws.Range("I70").Value = 4569.5884
$ws.Range("J70").Value = 4847.4375
$ws.Range("K70").Value = 4569.5884
$ws.Range("L70").Value = 4847.4375
$ws.Range("M70").Value = -4299.5884
$ws.Range("N70").Value = -5387.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4704.303
$ws.Range("I73").Value = 4569.5884
$ws.Range("J73").Value = 4847.4375
$ws.Range("K73").Value = 4569.5884
$ws.Range("L73").Value = 4847.4375
$ws.Range("M73").Value = -3633.5884
$ws.Range("N73").Value = -6719.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1533
$ws.Range("I113").Value = 1599.5
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1599.5
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 570.5
$ws.Range("N113").Value = -5740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1795.6
$ws.Range("I126").Value = 1541.2
$ws.Range("J126").Value = 2050
$ws.Range("K126").Value = 4623.6
$ws.Range("L126").Value = 6150
$ws.Range("M126").Value = -2153.6
$ws.Range("N126").Value = -11090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 439.6
$ws.Range("J22").Value = 474.5
$ws.Range("L22").Value = 474.5
$ws.Range("N22").Value = -1064.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 439.6
$ws.Range("J27").Value = 474.5
$ws.Range("L27").Value = 474.5
$ws.Range("N27").Value = -688.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2278.0952
$ws.Range("I40").Value = 2013.75
$ws.Range("J40").Value = 2630.5557
$ws.Range("K40").Value = 2013.75
$ws.Range("L40").Value = 2630.5557
$ws.Range("M40").Value = -1877.75
$ws.Range("N40").Value = -2902.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2254.5557
$ws.Range("I122").Value = 2050
$ws.Range("J122").Value = 2280.125
$ws.Range("K122").Value = 6150
$ws.Range("L122").Value = 6840.375
$ws.Range("M122").Value = -3700
$ws.Range("N122").Value = -11740.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2027.3846
$ws.Range("I122").Value = 2106.2222
$ws.Range("K122").Value = 6318.6666
$ws.Range("M122").Value = -3868.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3143.625
$ws.Range("J126").Value = 3973.75
$ws.Range("L126").Value = 11921.25
$ws.Range("N126").Value = -16861.25
